$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.713.80"
$ws.Range("D3").Value = "1.601.26"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "211.69"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.827.18"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.641.57"
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D17").Value = "26.695.44"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.10"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "144.17"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "1.293.48"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "0.601"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("D38").Value = "1.14"
$ws.Range("E38").Value = "  +6.73%  "
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "63.01"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").Value = "1.739.25"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "90.79"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").Value = "0.101"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -0.09%  "
